$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for all existing data rows
# (rows 2 through 205) from serial 45192 (2023-09-23) to serial 45202 (2023-10-03).
$ws.Range("C2:C205").Value = 45202

# Give row 205 an explicit row height (matches the rest of the data rows).
$ws.Rows.Item(205).RowHeight = 15

# Append the new record as row 206.
$ws.Range("A206").Value = "A 46444-2023"
$ws.Range("B206").Value = 45197
$ws.Range("C206").Value = 45202
$ws.Range("D206").Value = "ÖREBRO LÄN"
$ws.Range("E206").Value = "NORA"
$ws.Range("G206").Value = 1.1
$ws.Range("H206:Q206").Value = 0

# Match formatting used by the other rows: date format on B/C, wrapped text on R.
$ws.Range("B206:C206").NumberFormat = "YYYY-MM-DD"
$ws.Range("R206").WrapText = $true
